$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer (first page), currently holds id="3" Pearson logo: image1.png -> image2.png
$ftFirst = $sec.Footers.Item(2)
if ($ftFirst.Exists) {
    $shapes = $ftFirst.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $converted = $shp.ConvertToShape()
            $converted.Name = "image2.png"
            $converted.ConvertToInlineShape() | Out-Null
        }
    }
}

# --- Footer (default/primary), currently holds id="2" Pearson logo: image1.png -> image2.png
$ftDefault = $sec.Footers.Item(1)
if ($ftDefault.Exists) {
    $shapes = $ftDefault.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $converted = $shp.ConvertToShape()
            $converted.Name = "image2.png"
            $converted.ConvertToInlineShape() | Out-Null
        }
    }
}

# --- Header (first page), currently holds id="1" BTec logo: image2.jpg -> image1.jpg
$hdFirst = $sec.Headers.Item(2)
if ($hdFirst.Exists) {
    $shapes = $hdFirst.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $converted = $shp.ConvertToShape()
            $converted.Name = "image1.jpg"
            $converted.ConvertToInlineShape() | Out-Null
        }
    }
}

Write-Output "done"
